$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 0.5
$ws.Range("I2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.6666666666666666
$ws.Range("Q2").Value = 0.2

# Row 3
$ws.Range("F3").Value = 0.75
$ws.Range("I3").Value = 0.25
$ws.Range("M3").Value = 0.5
$ws.Range("Q3").Value = 0.25

# Row 4
$ws.Range("F4").Value = 0.6
$ws.Range("I4").Value = 0.2857142857142858
$ws.Range("M4").Value = 0.5714285714285715
$ws.Range("Q4").Value = 0.2222222222222222

# Row 5
$ws.Range("F5").Value = 0.6818181818181818
$ws.Range("I5").Value = 0.2631578947368421
$ws.Range("M5").Value = 0.5263157894736842
$ws.Range("Q5").Value = 0.2380952380952381

# Row 6
$ws.Range("F6").Value = 0.9552917361318499
$ws.Range("I6").Value = 0.3279120902703183
$ws.Range("M6").Value = 0.7712467209289329
$ws.Range("Q6").Value = 0.1412240503882725
